# Atualizacao de bases das ligas, do dia: 30-03-2024 as 19:32
#
# Two Singapore Premier League fixtures each had their full set of match
# data (ids, teams, score, odds) swapped with the data of the "sibling"
# row directly below/above it; the match number in column A stays put on
# its own row. Re-assert the final, corrected values for every cell of
# the eight affected rows (4,5 / 22,23 / 38,39 / 54,55).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value2 = 2
$ws.Range("B4").Value2 = 6228585
$ws.Range("C4").Value2 = "Singapore Premier League"
$ws.Range("D4").Value2 = "Singapore Premier League"
$ws.Range("E4").Value2 = 45083.36458333334
$ws.Range("F4").Value2 = "Albirex Niigata Singapore"
$ws.Range("G4").Value2 = "Geylang International"
$ws.Range("H4").Value2 = 3
$ws.Range("I4").Value2 = 0
$ws.Range("J4").Value2 = "H"
$ws.Range("K4").Value2 = 1.222
$ws.Range("L4").Value2 = 5.75
$ws.Range("M4").Value2 = 8
$ws.Range("N4").Value2 = 1.25
$ws.Range("O4").Value2 = 5.5
$ws.Range("P4").Value2 = 7
$ws.Range("Q4").Value2 = -1.75
$ws.Range("R4").Value2 = 1.875
$ws.Range("S4").Value2 = 1.975
$ws.Range("T4").Value2 = 4
$ws.Range("U4").Value2 = 1.825
$ws.Range("V4").Value2 = 2.025
$ws.Range("W4").Value2 = 0.25
$ws.Range("X4").Value2 = -1
$ws.Range("Y4").Value2 = -1
$ws.Range("Z4").Value2 = 0.875
$ws.Range("AA4").Value2 = -1
$ws.Range("AB4").Value2 = -1
$ws.Range("AC4").Value2 = 1.025

# Row 5
$ws.Range("A5").Value2 = 3
$ws.Range("B5").Value2 = 6228586
$ws.Range("C5").Value2 = "Singapore Premier League"
$ws.Range("D5").Value2 = "Singapore Premier League"
$ws.Range("E5").Value2 = 45083.36458333334
$ws.Range("F5").Value2 = "Young Lions"
$ws.Range("G5").Value2 = "Tanjong Pagar United"
$ws.Range("H5").Value2 = 3
$ws.Range("I5").Value2 = 4
$ws.Range("J5").Value2 = "A"
$ws.Range("K5").Value2 = 3.4
$ws.Range("L5").Value2 = 4.1
$ws.Range("M5").Value2 = 1.727
$ws.Range("N5").Value2 = 4.5
$ws.Range("O5").Value2 = 4.5
$ws.Range("P5").Value2 = 1.533
$ws.Range("Q5").Value2 = 1
$ws.Range("R5").Value2 = 2
$ws.Range("S5").Value2 = 1.85
$ws.Range("T5").Value2 = 3.5
$ws.Range("U5").Value2 = 1.925
$ws.Range("V5").Value2 = 1.925
$ws.Range("W5").Value2 = -1
$ws.Range("X5").Value2 = -1
$ws.Range("Y5").Value2 = 0.5329999999999999
$ws.Range("Z5").Value2 = 0
$ws.Range("AA5").Value2 = -0
$ws.Range("AB5").Value2 = 0.925
$ws.Range("AC5").Value2 = -1

# Row 22
$ws.Range("A22").Value2 = 20
$ws.Range("B22").Value2 = 6228600
$ws.Range("C22").Value2 = "Singapore Premier League"
$ws.Range("D22").Value2 = "Singapore Premier League"
$ws.Range("E22").Value2 = 45113.36458333334
$ws.Range("F22").Value2 = "Tanjong Pagar United"
$ws.Range("G22").Value2 = "Balestier Khalsa FC"
$ws.Range("H22").Value2 = 2
$ws.Range("I22").Value2 = 3
$ws.Range("J22").Value2 = "A"
$ws.Range("K22").Value2 = 3.2
$ws.Range("L22").Value2 = 4
$ws.Range("M22").Value2 = 1.8
$ws.Range("N22").Value2 = 3.4
$ws.Range("O22").Value2 = 4.2
$ws.Range("P22").Value2 = 1.8
$ws.Range("Q22").Value2 = 0.75
$ws.Range("R22").Value2 = 1.825
$ws.Range("S22").Value2 = 2.025
$ws.Range("T22").Value2 = 4.5
$ws.Range("U22").Value2 = 2
$ws.Range("V22").Value2 = 1.85
$ws.Range("W22").Value2 = -1
$ws.Range("X22").Value2 = -1
$ws.Range("Y22").Value2 = 0.8
$ws.Range("Z22").Value2 = -0.5
$ws.Range("AA22").Value2 = 0.5125
$ws.Range("AB22").Value2 = 1
$ws.Range("AC22").Value2 = -1

# Row 23
$ws.Range("A23").Value2 = 21
$ws.Range("B23").Value2 = 6228599
$ws.Range("C23").Value2 = "Singapore Premier League"
$ws.Range("D23").Value2 = "Singapore Premier League"
$ws.Range("E23").Value2 = 45113.36458333334
$ws.Range("F23").Value2 = "Hougang United FC"
$ws.Range("G23").Value2 = "Tampines Rovers FC"
$ws.Range("H23").Value2 = 0
$ws.Range("I23").Value2 = 1
$ws.Range("J23").Value2 = "A"
$ws.Range("K23").Value2 = 4.75
$ws.Range("L23").Value2 = 4.2
$ws.Range("M23").Value2 = 1.5
$ws.Range("N23").Value2 = 7.5
$ws.Range("O23").Value2 = 4.75
$ws.Range("P23").Value2 = 1.3
$ws.Range("Q23").Value2 = 1.5
$ws.Range("R23").Value2 = 2
$ws.Range("S23").Value2 = 1.85
$ws.Range("T23").Value2 = 3.75
$ws.Range("U23").Value2 = 2
$ws.Range("V23").Value2 = 1.85
$ws.Range("W23").Value2 = -1
$ws.Range("X23").Value2 = -1
$ws.Range("Y23").Value2 = 0.3
$ws.Range("Z23").Value2 = 1
$ws.Range("AA23").Value2 = -1
$ws.Range("AB23").Value2 = -1
$ws.Range("AC23").Value2 = 0.8500000000000001

# Row 38
$ws.Range("A38").Value2 = 36
$ws.Range("B38").Value2 = 6228611
$ws.Range("C38").Value2 = "Singapore Premier League"
$ws.Range("D38").Value2 = "Singapore Premier League"
$ws.Range("E38").Value2 = 45135.36458333334
$ws.Range("F38").Value2 = "Albirex Niigata Singapore"
$ws.Range("G38").Value2 = "Tampines Rovers FC"
$ws.Range("H38").Value2 = 6
$ws.Range("I38").Value2 = 3
$ws.Range("J38").Value2 = "H"
$ws.Range("K38").Value2 = 1.7
$ws.Range("L38").Value2 = 4
$ws.Range("M38").Value2 = 3.6
$ws.Range("N38").Value2 = 1.5
$ws.Range("O38").Value2 = 3.8
$ws.Range("P38").Value2 = 5.25
$ws.Range("Q38").Value2 = -1.25
$ws.Range("R38").Value2 = 1.975
$ws.Range("S38").Value2 = 1.875
$ws.Range("T38").Value2 = 4
$ws.Range("U38").Value2 = 2.025
$ws.Range("V38").Value2 = 1.825
$ws.Range("W38").Value2 = 0.5
$ws.Range("X38").Value2 = -1
$ws.Range("Y38").Value2 = -1
$ws.Range("Z38").Value2 = 0.9750000000000001
$ws.Range("AA38").Value2 = -1
$ws.Range("AB38").Value2 = 1.025
$ws.Range("AC38").Value2 = -1

# Row 39
$ws.Range("A39").Value2 = 37
$ws.Range("B39").Value2 = 6228613
$ws.Range("C39").Value2 = "Singapore Premier League"
$ws.Range("D39").Value2 = "Singapore Premier League"
$ws.Range("E39").Value2 = 45135.36458333334
$ws.Range("F39").Value2 = "DPMM FC"
$ws.Range("G39").Value2 = "Geylang International"
$ws.Range("H39").Value2 = 1
$ws.Range("I39").Value2 = 2
$ws.Range("J39").Value2 = "A"
$ws.Range("K39").Value2 = 2.5
$ws.Range("L39").Value2 = 3.75
$ws.Range("M39").Value2 = 2.25
$ws.Range("N39").Value2 = 2.45
$ws.Range("O39").Value2 = 3.6
$ws.Range("P39").Value2 = 2.3
$ws.Range("Q39").Value2 = 0
$ws.Range("R39").Value2 = 2
$ws.Range("S39").Value2 = 1.85
$ws.Range("T39").Value2 = 3.75
$ws.Range("U39").Value2 = 1.925
$ws.Range("V39").Value2 = 1.925
$ws.Range("W39").Value2 = -1
$ws.Range("X39").Value2 = -1
$ws.Range("Y39").Value2 = 1.3
$ws.Range("Z39").Value2 = -1
$ws.Range("AA39").Value2 = 0.8500000000000001
$ws.Range("AB39").Value2 = -1
$ws.Range("AC39").Value2 = 0.925

# Row 54
$ws.Range("A54").Value2 = 52
$ws.Range("B54").Value2 = 7098763
$ws.Range("C54").Value2 = "Singapore Premier League"
$ws.Range("D54").Value2 = "Singapore Premier League"
$ws.Range("E54").Value2 = 45184.36458333334
$ws.Range("F54").Value2 = "Balestier Khalsa FC"
$ws.Range("G54").Value2 = "Tampines Rovers FC"
$ws.Range("H54").Value2 = 1
$ws.Range("I54").Value2 = 3
$ws.Range("J54").Value2 = "A"
$ws.Range("K54").Value2 = 5.25
$ws.Range("L54").Value2 = 4.2
$ws.Range("M54").Value2 = 1.5
$ws.Range("N54").Value2 = 5
$ws.Range("O54").Value2 = 4.5
$ws.Range("P54").Value2 = 1.45
$ws.Range("Q54").Value2 = 1.25
$ws.Range("R54").Value2 = 2
$ws.Range("S54").Value2 = 1.85
$ws.Range("T54").Value2 = 5
$ws.Range("U54").Value2 = 1.925
$ws.Range("V54").Value2 = 1.925
$ws.Range("W54").Value2 = -1
$ws.Range("X54").Value2 = -1
$ws.Range("Y54").Value2 = 0.45
$ws.Range("Z54").Value2 = -1
$ws.Range("AA54").Value2 = 0.8500000000000001
$ws.Range("AB54").Value2 = -1
$ws.Range("AC54").Value2 = 0.925

# Row 55
$ws.Range("A55").Value2 = 53
$ws.Range("B55").Value2 = 7094656
$ws.Range("C55").Value2 = "Singapore Premier League"
$ws.Range("D55").Value2 = "Singapore Premier League"
$ws.Range("E55").Value2 = 45184.36458333334
$ws.Range("F55").Value2 = "Tanjong Pagar United"
$ws.Range("G55").Value2 = "DPMM FC"
$ws.Range("H55").Value2 = 1
$ws.Range("I55").Value2 = 1
$ws.Range("J55").Value2 = "D"
$ws.Range("K55").Value2 = 2.15
$ws.Range("L55").Value2 = 3.75
$ws.Range("M55").Value2 = 2.7
$ws.Range("N55").Value2 = 2.1
$ws.Range("O55").Value2 = 4.2
$ws.Range("P55").Value2 = 2.625
$ws.Range("Q55").Value2 = -0.25
$ws.Range("R55").Value2 = 1.925
$ws.Range("S55").Value2 = 1.925
$ws.Range("T55").Value2 = 4.25
$ws.Range("U55").Value2 = 1.9
$ws.Range("V55").Value2 = 1.95
$ws.Range("W55").Value2 = -1
$ws.Range("X55").Value2 = 3.2
$ws.Range("Y55").Value2 = -1
$ws.Range("Z55").Value2 = -0.5
$ws.Range("AA55").Value2 = 0.4625
$ws.Range("AB55").Value2 = -1
$ws.Range("AC55").Value2 = 0.95
